$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the required fields in the header row.
$ws.Range("A1").Value = "Email (*Required)"
$ws.Range("B1").Value = "First name (*Required)"
$ws.Range("C1").Value = "Last name (*Required)"

# Widen the columns that now hold the longer header text.
$ws.Columns.Item(2).ColumnWidth = 19.833333333333332
$ws.Columns.Item(3).ColumnWidth = 21.333333333333332
$ws.Columns.Item(4).ColumnWidth = 19.5
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668

# Move the active selection.
$null = $ws.Range("C5").Select()
